$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4030
$ws.Range("C3").Value = 4030
$ws.Range("C4").Value = 4030
$ws.Range("C5").Value = 4390
$ws.Range("C6").Value = 4708
$ws.Range("C7").Value = 4708
$ws.Range("C8").Value = 4708
$ws.Range("C9").Value = 4831
$ws.Range("C10").Value = 5248
$ws.Range("C11").Value = 5291
$ws.Range("C12").Value = 5304
$ws.Range("C13").Value = 5339
$ws.Range("C14").Value = 5502
$ws.Range("C15").Value = 5502
